# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets
# to reflect the latest scrape, as described by the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 7225
$wsExpo.Range("F7").Value = 184
$wsExpo.Range("F8").Value = 131
$wsExpo.Range("F12").Value = 219
$wsExpo.Range("F13").Value = 14
$wsExpo.Range("F14").Value = 461
$wsExpo.Range("F16").Value = 1857
$wsExpo.Range("F18").Value = 40
$wsExpo.Range("F19").Value = 3765
$wsExpo.Range("F25").Value = 35
$wsExpo.Range("F26").Value = 2433
$wsExpo.Range("F28").Value = 300
$wsExpo.Range("F33").Value = 17
$wsExpo.Range("F34").Value = 3
$wsExpo.Range("F38").Value = 1456
$wsExpo.Range("F39").Value = 152

# Sheet "全部类型" (all types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 7225
$wsAll.Range("F8").Value = 184
$wsAll.Range("F9").Value = 131
$wsAll.Range("F13").Value = 219
$wsAll.Range("F14").Value = 14
$wsAll.Range("F15").Value = 461
$wsAll.Range("F17").Value = 1857
$wsAll.Range("F19").Value = 40
$wsAll.Range("F20").Value = 3765
$wsAll.Range("F26").Value = 35
$wsAll.Range("F27").Value = 2433
$wsAll.Range("F29").Value = 300
$wsAll.Range("F34").Value = 17
$wsAll.Range("F35").Value = 3
$wsAll.Range("F39").Value = 1456
$wsAll.Range("F40").Value = 152
